# Apply the block-order reshuffle to the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: header labels (columns C:F get reassigned) ---
$ws.Range("C1").Value = "bedrooms_1"
$ws.Range("D1").Value = "kitchens_2"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "living_rooms_2"

# --- Rows 2-7: the 0/1 indicator grid (columns A:F) ---
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 0
$ws.Cells.Item(2, 3).Value = 0
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0

$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0

$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 0

$ws.Cells.Item(5, 1).Value = 0
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 0

$ws.Cells.Item(6, 1).Value = 1
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 0
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 0

$ws.Cells.Item(7, 1).Value = 0
$ws.Cells.Item(7, 2).Value = 0
$ws.Cells.Item(7, 3).Value = 0
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 1
